$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 117 and 118: swap all data columns (B:AD), keep column A (id) fixed per row
$r117 = $ws.Range("B117:AD117").Value2
$r118 = $ws.Range("B118:AD118").Value2

$ws.Range("B117:AD117").Value2 = $r118
$ws.Range("B118:AD118").Value2 = $r117

# Rows 234, 236, 237: cyclic rotation of all data columns (B:AD), keep column A (id) fixed per row
# new234 = old236 ; new236 = old237 ; new237 = old234
$r234 = $ws.Range("B234:AD234").Value2
$r236 = $ws.Range("B236:AD236").Value2
$r237 = $ws.Range("B237:AD237").Value2

$ws.Range("B234:AD234").Value2 = $r236
$ws.Range("B236:AD236").Value2 = $r237
$ws.Range("B237:AD237").Value2 = $r234

Write-Host "done"
